$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 2913.066666666667, 3201, 2602, 0.05184965133666992),
    @(1, 2815.3, 3117, 2493, 0.05181701183319092),
    @(2, 2930.6, 3208, 2688, 0.05476129055023193),
    @(3, 3041.4, 3301, 2665, 0.05440890789031982),
    @(4, 2927.766666666667, 3187, 2708, 0.05468142032623291),
    @(5, 2927.9, 3247, 2527, 0.05462076663970947),
    @(6, 3279.933333333333, 3587, 3063, 0.05448505878448486),
    @(7, 2692.466666666667, 2913, 2383, 0.05223522981007894),
    @(8, 3034.5, 3322, 2496, 0.05160330931345622),
    @(9, 2836.7, 3081, 2449, 0.05213772455851237)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
